# fix: register & update: create & edit merchant
#
# Inserts a new "Merchant Type" header column between "Merchant Name" and
# "Merchant Email" on the header row (row 3), shifting every header from
# the old column C onward one column to the right (C..Q -> D..R), without
# touching the pre-formatted (but empty) data-row cells in row 4, which
# keep their original column letters / styles. Also nudges the active
# selection to A4 and forces portrait page orientation, matching the
# saved workbook's new state.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column R (18) is brand new for row 3 -- give it the same header
# formatting (bold/fill/border/center, style index 1) as the rest of the
# header row before we drop a value into it.
$ws.Range("Q3").Copy()
$ws.Range("R3").PasteSpecial(-4122)  # xlPasteFormats

# Capture the existing header row (row 3) text, C3 through Q3, before
# overwriting anything -- we need the old values to shift them right.
$oldHeaders = @()
for ($col = 3; $col -le 17; $col++) {
    $oldHeaders += $ws.Cells.Item(3, $col).Value2
}

# Write the shifted headers back starting one column further right
# (D..R), walking from the right-most column towards C so we never
# clobber a value we still need to read (defensive; values were already
# captured above, but keep the safe order anyway).
for ($i = $oldHeaders.Length - 1; $i -ge 0; $i--) {
    $destCol = 4 + $i
    $ws.Cells.Item(3, $destCol).Value = $oldHeaders[$i]
}

# New header inserted at column C.
$ws.Range("C3").Value = "Merchant Type"

# Move the active selection to A4, as recorded in the saved view state.
[void]$ws.Range("A4").Select()

# Force explicit portrait orientation (adds <pageSetup .../> on save).
$ws.PageSetup.Orientation = 1
